$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill previously empty Spanish (column C) translation cells
$ws.Range("C45").Value = "Jetpack"
$ws.Range("C46").Value = "Flecha"
$ws.Range("C47").Value = "Taser"
$ws.Range("C57").Value = "Al Recoger"
$ws.Range("C125").Value = "¿Rotar?"
$ws.Range("C126").Value = "Velocidad de Rotación"
$ws.Range("C127").Value = "Munición"
$ws.Range("C128").Value = "Taser infinito"

# Row 54 no longer needs its custom (larger) height now the cell is filled
$ws.Rows.Item(54).AutoFit()

# Update the active selection to the next empty row (C129)
[void]$ws.Range("C129").Select()
